# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
# 展览: rows 2-10 map directly to F2..F10
# 全部类型: same events but with one extra row (row 4) inserted in the
#           middle, so the matching rows are F2, F3, F5..F11.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 722
$ws1.Range("F3").Value = 37
$ws1.Range("F4").Value = 242
$ws1.Range("F5").Value = 2639
$ws1.Range("F6").Value = 55
$ws1.Range("F7").Value = 3674
$ws1.Range("F8").Value = 467
$ws1.Range("F9").Value = 929
$ws1.Range("F10").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 722
$ws4.Range("F3").Value = 37
$ws4.Range("F5").Value = 242
$ws4.Range("F6").Value = 2639
$ws4.Range("F7").Value = 55
$ws4.Range("F8").Value = 3674
$ws4.Range("F9").Value = 467
$ws4.Range("F10").Value = 929
$ws4.Range("F11").Value = 3
